$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.054.29'
$ws.Range('E2').Value = '  +1.13%  '

$ws.Range('D3').Value = '2.524.75'
$ws.Range('E3').Value = '  -0.94%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = "''596.06"
$ws.Range('E5').Value = '  +0.91%  '

$ws.Range('D6').Value = "''174.55"
$ws.Range('E6').Value = '  -0.08%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('D9').Value = '2.524.30'
$ws.Range('E9').Value = '  -0.92%  '

$ws.Range('D10').Value = "''0.140"
$ws.Range('E10').Value = '  +0.52%  '

$ws.Range('E11').Value = '  +2.11%  '

$ws.Range('E12').Value = '  -1.01%  '

$ws.Range('E13').Value = '  -2.13%  '

$ws.Range('D14').Value = "''26.59"
$ws.Range('E14').Value = '  -1.64%  '

$ws.Range('D15').Value = '2.985.83'
$ws.Range('E15').Value = '  -0.91%  '

$ws.Range('D16').Value = "''0.0000178"
$ws.Range('E16').Value = '  +0.06%  '

$ws.Range('D17').Value = '68.015.94'
$ws.Range('E17').Value = '  +1.25%  '

$ws.Range('D18').Value = '2.503.37'
$ws.Range('E18').Value = '  -1.52%  '

$ws.Range('D19').Value = "''11.95"
$ws.Range('E19').Value = '  +4.88%  '

$ws.Range('D20').Value = "''8.11"
$ws.Range('E20').Value = '  +0.59%  '

$ws.Range('D21').Value = "''364.43"
$ws.Range('E21').Value = '  +2.60%  '

$ws.Range('E22').Value = '  -1.33%  '

$ws.Range('E23').Value = '  -0.26%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = "''1.00"
$ws.Range('E24').Value = '  +0.19%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = "''71.28"
$ws.Range('E25').Value = '  +1.67%  '

$ws.Range('E26').Value = '  -3.66%  '

$ws.Range('D27').Value = "''10.03"
$ws.Range('E27').Value = '  +0.28%  '

$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.650.84'
$ws.Range('E28').Value = '  -1.28%  '

$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = "''0.993"
$ws.Range('E29').Value = '  -0.80%  '

$ws.Range('D30').Value = '0.0₃0984'
$ws.Range('E30').Value = '  -1.31%  '

$ws.Range('D31').Value = "''8.36"
$ws.Range('E31').Value = '  +1.35%  '

$ws.Range('D32').Value = "''533.43"
$ws.Range('E32').Value = '  -0.49%  '

$ws.Range('E33').Value = '  +2.44%  '

$ws.Range('E34').Value = '  -1.51%  '

$ws.Range('E35').Value = '  -2.54%  '

$ws.Range('E36').Value = '  +0.00%  '

$ws.Range('D37').Value = "''157.88"
$ws.Range('E37').Value = '  +0.21%  '

$ws.Range('E38').Value = '  -2.16%  '

$ws.Range('D39').Value = "''18.76"
$ws.Range('E39').Value = '  -0.27%  '

$ws.Range('D40').Value = "''18.67"
$ws.Range('E40').Value = '  +1.27%  '

$ws.Range('E41').Value = '  -0.60%  '

$ws.Range('E42').Value = '  -0.33%  '

$ws.Range('E43').Value = '  -1.62%  '

$ws.Range('E44').Value = '  -1.94%  '

$ws.Range('E45').Value = '  +0.01%  '

$ws.Range('E46').Value = '  -2.63%  '

$ws.Range('D47').Value = '0.0⁦0279'
$ws.Range('E47').Value = '  -0.45%  '

$ws.Range('D48').Value = "''0.556"
$ws.Range('E48').Value = '  -1.59%  '

$ws.Range('E49').Value = '  -0.09%  '

$ws.Range('E50').Value = '  +0.62%  '

$ws.Range('E51').Value = '  -1.28%  '
